# "Generate Report for Handback" - add a new handback-failure row
# (a56b4c39-1622-461c-be84-e126b5128073) ahead of the trailing
# ".localization-config" row on every sheet.
#
# NOTE: cell values are written in the exact top-to-bottom / left-to-right,
# sheet-by-sheet order that the final workbook lays them out in, so that
# newly-introduced strings land in the shared-string table in the same
# order the original report generator produced.

$wb = $excel.ActiveWorkbook

$newMdName   = "a56b4c39-1622-461c-be84-e126b5128073.md"
$newZhXlf    = "a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.zh-cn.xlf"
$newDeXlf    = "a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.de-de.xlf"
$failedText  = "Handback transform failed"
$zhDateTime  = "2016-03-09 03:18:30"
$deDateTime  = "2016-03-09 03:18:41"
$epoch       = "0001-01-01 00:00:00"
$includeTxt  = "Include"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Make room for the new row, shifting the old row 3 (".localization-config")
# down to row 4. Whole-row insert copies formatting from the row above,
# which happens to already match what the new row needs.
$ws1.Rows.Item(3).Insert()

# Re-establish rows 1-2 in-order first (no-ops content-wise, but keeps the
# shared-string discovery order identical to a from-scratch regeneration).
$ws1.Cells.Item(1,1).Value2 = "File Name"
$ws1.Cells.Item(1,2).Value2 = "zh-cn"
$ws1.Cells.Item(1,3).Value2 = "de-de"
$ws1.Cells.Item(2,1).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$ws1.Cells.Item(2,2).Value2 = "Handed back: in sync with en-US"
$ws1.Cells.Item(2,3).Value2 = "Handed back: in sync with en-US"

# New row 3.
$ws1.Cells.Item(3,1).Value2 = $newMdName
$ws1.Cells.Item(3,2).Value2 = $failedText
$ws1.Cells.Item(3,3).Value2 = $failedText

# Row 4 already holds the shifted-down original row-3 values/styles.

# Rebuild hyperlinks for this sheet (row-shift does not move the anchor of
# an existing hyperlink, so clear and re-add them all in final position).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/a56b4c39-1622-461c-be84-e126b5128073.md", "", "", $newMdName)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Rows.Item(3).Insert()

$ws2.Cells.Item(1,1).Value2 = "Source File Name"
$ws2.Cells.Item(1,2).Value2 = "Status"
$ws2.Cells.Item(1,3).Value2 = "Latest Handoff File"
$ws2.Cells.Item(1,4).Value2 = "Latest Handoff Datetime"
$ws2.Cells.Item(1,5).Value2 = "Latest Target File"
$ws2.Cells.Item(1,6).Value2 = "Latest Handback File"
$ws2.Cells.Item(1,7).Value2 = "Latest Handback DateTime"
$ws2.Cells.Item(1,8).Value2 = "Handoff Reason"
$ws2.Cells.Item(1,9).Value2 = "Dependency From"

$ws2.Cells.Item(2,1).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$ws2.Cells.Item(2,2).Value2 = "Handed back: in sync with en-US"
$ws2.Cells.Item(2,3).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf"
$ws2.Cells.Item(2,4).Value2 = "2016-03-09 03:15:56"
$ws2.Cells.Item(2,5).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$ws2.Cells.Item(2,6).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf"
$ws2.Cells.Item(2,7).Value2 = "2016-03-09 03:16:53"
$ws2.Cells.Item(2,8).Value2 = $includeTxt

# New row 3.
$ws2.Cells.Item(3,1).Value2 = $newMdName
$ws2.Cells.Item(3,2).Value2 = $failedText
$ws2.Cells.Item(3,3).Value2 = $newZhXlf
$ws2.Cells.Item(3,4).Value2 = $zhDateTime
$ws2.Cells.Item(3,5).ClearFormats()
$ws2.Cells.Item(3,6).ClearFormats()
$ws2.Cells.Item(3,7).Value2 = $epoch
$ws2.Cells.Item(3,8).Value2 = $includeTxt

# Row 4 already holds the shifted-down original row-3 values/styles.

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3d55c7d0eb607ffa99726cb296d186886cfb9af/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3504c7f8f7f183c35f82b0409aa4baef8c0ec3a1/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e203a7323ea4c108dc217da020b3be848c464830/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/a56b4c39-1622-461c-be84-e126b5128073.md", "", "", $newMdName)
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e203a7323ea4c108dc217da020b3be848c464830/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.zh-cn.xlf", "", "", $newZhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Rows.Item(3).Insert()

$ws3.Cells.Item(1,1).Value2 = "Source File Name"
$ws3.Cells.Item(1,2).Value2 = "Status"
$ws3.Cells.Item(1,3).Value2 = "Latest Handoff File"
$ws3.Cells.Item(1,4).Value2 = "Latest Handoff Datetime"
$ws3.Cells.Item(1,5).Value2 = "Latest Target File"
$ws3.Cells.Item(1,6).Value2 = "Latest Handback File"
$ws3.Cells.Item(1,7).Value2 = "Latest Handback DateTime"
$ws3.Cells.Item(1,8).Value2 = "Handoff Reason"
$ws3.Cells.Item(1,9).Value2 = "Dependency From"

$ws3.Cells.Item(2,1).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$ws3.Cells.Item(2,2).Value2 = "Handed back: in sync with en-US"
$ws3.Cells.Item(2,3).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf"
$ws3.Cells.Item(2,4).Value2 = "2016-03-09 03:16:11"
$ws3.Cells.Item(2,5).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.md"
$ws3.Cells.Item(2,6).Value2 = "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf"
$ws3.Cells.Item(2,7).Value2 = "2016-03-09 03:17:29"
$ws3.Cells.Item(2,8).Value2 = $includeTxt

# New row 3.
$ws3.Cells.Item(3,1).Value2 = $newMdName
$ws3.Cells.Item(3,2).Value2 = $failedText
$ws3.Cells.Item(3,3).Value2 = $newDeXlf
$ws3.Cells.Item(3,4).Value2 = $deDateTime
$ws3.Cells.Item(3,5).ClearFormats()
$ws3.Cells.Item(3,6).ClearFormats()
$ws3.Cells.Item(3,7).Value2 = $epoch
$ws3.Cells.Item(3,8).Value2 = $includeTxt

# Row 4 already holds the shifted-down original row-3 values/styles.

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8aa48bd59e50644c8a69fd21ea7db36aca6b67a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5c22e97140677a61ceb971650cf23cd623a0eeea/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cd3e808463306f071d0b40f4a59fbbaa5d7f7ad0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/a56b4c39-1622-461c-be84-e126b5128073.md", "", "", $newMdName)
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cd3e808463306f071d0b40f4a59fbbaa5d7f7ad0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.de-de.xlf", "", "", $newDeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/.localization-config", "", "", ".localization-config")

Write-Output "done"
